$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value2 = 0.2080200501253133
$ws.Range("C2").Value2 = 0.543859649122807
$ws.Range("J2").Value2 = 0.01503759398496241
$ws.Range("P2").Value2 = 0.1428571428571428
$ws.Range("S2").Value2 = 0.09022556390977443
$ws.Range("B3").Value2 = 0.013215859030837
$ws.Range("C3").Value2 = 0.03524229074889868
$ws.Range("J3").Value2 = 0.03964757709251102
$ws.Range("P3").Value2 = 0.7048458149779736
$ws.Range("S3").Value2 = 0.2070484581497797
$ws.Range("J4").Value2 = 0.02325581395348837
$ws.Range("P4").Value2 = 0.6511627906976745
$ws.Range("S4").Value2 = 0.3255813953488372
$ws.Range("B6").Value2 = 0.07207207207207207
$ws.Range("D6").Value2 = 0.01351351351351351
$ws.Range("F6").Value2 = 0.05405405405405406
$ws.Range("J6").Value2 = 0.2342342342342342
$ws.Range("O6").Value2 = 0.01801801801801802
$ws.Range("Q6").Value2 = 0.1396396396396396
$ws.Range("R6").Value2 = 0.09009009009009009
$ws.Range("S6").Value2 = 0.3783783783783784
$ws.Range("B7").Value2 = 0.09836065573770492
$ws.Range("D7").Value2 = 0.01092896174863388
$ws.Range("F7").Value2 = 0.03278688524590164
$ws.Range("J7").Value2 = 0.2021857923497268
$ws.Range("O7").Value2 = 0.01092896174863388
$ws.Range("Q7").Value2 = 0.185792349726776
$ws.Range("R7").Value2 = 0.06557377049180328
$ws.Range("S7").Value2 = 0.3934426229508197
$ws.Range("B8").Value2 = 0.1326781326781327
$ws.Range("D8").Value2 = 0.02702702702702703
$ws.Range("F8").Value2 = 0.09582309582309582
$ws.Range("J8").Value2 = 0.1154791154791155
$ws.Range("O8").Value2 = 0.01228501228501228
$ws.Range("Q8").Value2 = 0.1646191646191646
$ws.Range("R8").Value2 = 0.08353808353808354
$ws.Range("S8").Value2 = 0.3685503685503685
$ws.Range("B9").Value2 = 0.07142857142857142
$ws.Range("D9").Value2 = 0.00510204081632653
$ws.Range("E9").Value2 = 0.00510204081632653
$ws.Range("F9").Value2 = 0.06122448979591837
$ws.Range("J9").Value2 = 0.1479591836734694
$ws.Range("O9").Value2 = 0.02040816326530612
$ws.Range("Q9").Value2 = 0.1989795918367347
$ws.Range("R9").Value2 = 0.07653061224489796
$ws.Range("S9").Value2 = 0.413265306122449
$ws.Range("B10").Value2 = 0.1422623178348369
$ws.Range("D10").Value2 = 0.01804302567661346
$ws.Range("E10").Value2 = 0.002775850104094379
$ws.Range("F10").Value2 = 0.06453851492019431
$ws.Range("J10").Value2 = 0.1269951422623178
$ws.Range("O10").Value2 = 0.01665510062456627
$ws.Range("Q10").Value2 = 0.2095766828591256
$ws.Range("R10").Value2 = 0.08119361554476058
$ws.Range("S10").Value2 = 0.3379597501734906
$ws.Range("G11").Value2 = 0.1643835616438356
$ws.Range("J11").Value2 = 0.0821917808219178
$ws.Range("K11").Value2 = 0.1986301369863014
$ws.Range("L11").Value2 = 0.547945205479452
$ws.Range("S11").Value2 = 0.00684931506849315
$ws.Range("G12").Value2 = 0.6909090909090909
$ws.Range("J12").Value2 = 0.2666666666666667
$ws.Range("L12").Value2 = 0.01818181818181818
$ws.Range("S12").Value2 = 0.02424242424242424
$ws.Range("F13").Value2 = 0.02631578947368421
$ws.Range("G13").Value2 = 0.6842105263157895
$ws.Range("J13").Value2 = 0.2368421052631579
$ws.Range("S13").Value2 = 0.05263157894736842
$ws.Range("F15").Value2 = 0.01463414634146342
$ws.Range("H15").Value2 = 0.1268292682926829
$ws.Range("I15").Value2 = 0.08780487804878048
$ws.Range("J15").Value2 = 0.3902439024390244
$ws.Range("K15").Value2 = 0.04390243902439024
$ws.Range("M15").Value2 = 0.01951219512195122
$ws.Range("O15").Value2 = 0.03414634146341464
$ws.Range("S15").Value2 = 0.2829268292682927
$ws.Range("F16").Value2 = 0.02100840336134454
$ws.Range("H16").Value2 = 0.1302521008403361
$ws.Range("I16").Value2 = 0.08403361344537816
$ws.Range("J16").Value2 = 0.4495798319327731
$ws.Range("K16").Value2 = 0.09663865546218488
$ws.Range("M16").Value2 = 0.02100840336134454
$ws.Range("O16").Value2 = 0.02100840336134454
$ws.Range("S16").Value2 = 0.1764705882352941
$ws.Range("F17").Value2 = 0.002127659574468085
$ws.Range("H17").Value2 = 0.1659574468085106
$ws.Range("I17").Value2 = 0.07659574468085106
$ws.Range("J17").Value2 = 0.4723404255319149
$ws.Range("K17").Value2 = 0.08723404255319149
$ws.Range("M17").Value2 = 0.01276595744680851
$ws.Range("O17").Value2 = 0.0425531914893617
$ws.Range("S17").Value2 = 0.1404255319148936
$ws.Range("F18").Value2 = 0.01020408163265306
$ws.Range("H18").Value2 = 0.1836734693877551
$ws.Range("I18").Value2 = 0.1071428571428571
$ws.Range("J18").Value2 = 0.3979591836734694
$ws.Range("K18").Value2 = 0.08163265306122448
$ws.Range("M18").Value2 = 0.01530612244897959
$ws.Range("O18").Value2 = 0.08673469387755102
$ws.Range("S18").Value2 = 0.1173469387755102
$ws.Range("F19").Value2 = 0.01469450889404486
$ws.Range("H19").Value2 = 0.1832946635730859
$ws.Range("I19").Value2 = 0.07965970610982212
$ws.Range("J19").Value2 = 0.4075792730085073
$ws.Range("K19").Value2 = 0.1098221191028616
$ws.Range("M19").Value2 = 0.01701469450889405
$ws.Range("N19").Value2 = 0.0007733952049497294
$ws.Range("O19").Value2 = 0.0711523588553751
$ws.Range("S19").Value2 = 0.1160092807424594
